$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(23).RowHeight = 49.2
$ws.Range("D23").Value = "→         →"
$ws.Range("E23").Value = "1 K ohm Resistor, 1 K ohm resistor"
$ws.Range("F23").Value = "→                →"
$ws.Range("G23").Value = "UNO Reset"
$ws.Range("H23").Value = ", Ground"
$ws.Range("D23").Copy()
$ws.Range("F23:H23").PasteSpecial(-4122)
